# Apply cell content updates to Sheet1 per commit diff (cryptos list refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds text-formatted price strings (e.g. "28.552.22", "1.151").
# Force text format across the data range first so Excel does not
# auto-convert the numeric-looking replacement strings into numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '28.552.22'
$ws.Range("E2").Value = '  +1.28%  '

$ws.Range("D3").Value = '1.912.12'
$ws.Range("E3").Value = '  +4.60%  '

$ws.Range("E4").Value = '  +0.11%  '

$ws.Range("D5").Value = '315.24'
$ws.Range("E5").Value = '  +1.37%  '

$ws.Range("E6").Value = '  +0.03%  '

$ws.Range("E7").Value = '  +3.98%  '

$ws.Range("E8").Value = '  +0.54%  '

$ws.Range("D9").Value = '0.09680'
$ws.Range("E9").Value = '  -2.41%  '

$ws.Range("D10").Value = '1.151'
$ws.Range("E10").Value = '  +3.73%  '

$ws.Range("D11").Value = '42.14'
$ws.Range("E11").Value = '  +1.93%  '

$ws.Range("D12").Value = '6.527'
$ws.Range("E12").Value = '  +1.25%  '

$ws.Range("D13").Value = '21.22'
$ws.Range("E13").Value = '  +2.76%  '

$ws.Range("D14").Value = '1.914.34'
$ws.Range("E14").Value = '  +4.63%  '

$ws.Range("D15").Value = '7.500'
$ws.Range("E15").Value = '  +2.59%  '

$ws.Range("E16").Value = '  +0.07%  '

$ws.Range("D17").Value = '94.64'
$ws.Range("E17").Value = '  +1.86%  '

$ws.Range("D18").Value = '0.00001134'
$ws.Range("E18").Value = '  -0.90%  '

$ws.Range("D19").Value = '0.06655'
$ws.Range("E19").Value = '  +0.05%  '

$ws.Range("E20").Value = '  +5.65%  '

$ws.Range("D21").Value = '1.000'
$ws.Range("E21").Value = '  -0.01%  '

$ws.Range("D22").Value = '6.308'
$ws.Range("E22").Value = '  +5.10%  '

$ws.Range("D23").Value = '28.611.02'
$ws.Range("E23").Value = '  +1.34%  '

$ws.Range("D24").Value = '11.50'
$ws.Range("E24").Value = '  +1.31%  '

$ws.Range("D25").Value = '2.315'
$ws.Range("E25").Value = '  +3.18%  '

$ws.Range("D26").Value = '2.688'
$ws.Range("E26").Value = '  +10.90%  '

$ws.Range("D27").Value = '2.134.25'
$ws.Range("E27").Value = '  +4.66%  '

$ws.Range("D28").Value = '21.25'
$ws.Range("E28").Value = '  +2.37%  '

$ws.Range("D29").Value = '158.01'

$ws.Range("E30").Value = '  +1.38%  '

$ws.Range("D31").Value = '1.109'
$ws.Range("E31").Value = '  +6.42%  '

$ws.Range("E32").Value = '  +2.38%  '

$ws.Range("E33").Value = '  +2.90%  '

$ws.Range("D34").Value = '3.632'
$ws.Range("E34").Value = '  +0.99%  '

$ws.Range("D35").Value = '9.982'
$ws.Range("E35").Value = '  +10.23%  '

$ws.Range("D36").Value = '0.06793'
$ws.Range("E36").Value = '  +0.55%  '

$ws.Range("D37").Value = '1.280'
$ws.Range("E37").Value = '  +8.26%  '

$ws.Range("D38").Value = '0.02434'
$ws.Range("E38").Value = '  +3.79%  '

$ws.Range("D39").Value = '0.2216'
$ws.Range("E39").Value = '  +2.89%  '

$ws.Range("D40").Value = '11.81'
$ws.Range("E40").Value = '  +3.58%  '

$ws.Range("B41").Value = 'InternetComputer(DFINITY)'
$ws.Range("C41").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D41").Value = '5.090'
$ws.Range("E41").Value = '  +2.17%  '

$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D42").Value = '0.6466'
$ws.Range("E42").Value = '  +3.88%  '

$ws.Range("D43").Value = '1.191'
$ws.Range("E43").Value = '  +0.80%  '

$ws.Range("D44").Value = '1.000'
$ws.Range("E44").Value = '  +0.03%  '

$ws.Range("D45").Value = '13.55'
$ws.Range("E45").Value = '  +2.94%  '

$ws.Range("D46").Value = '0.6099'
$ws.Range("E46").Value = '  +2.64%  '

$ws.Range("D47").Value = '3.778'
$ws.Range("E47").Value = '  +2.15%  '

$ws.Range("D48").Value = '1.283'
$ws.Range("E48").Value = '  +0.89%  '

$ws.Range("E49").Value = '  +4.46%  '

$ws.Range("D50").Value = '125.10'
$ws.Range("E50").Value = '  +0.71%  '

$ws.Range("D51").Value = '1.204'
$ws.Range("E51").Value = '  +1.71%  '
